$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the result cell for the first row with "ok"
$ws.Range("B2").Value = "ok"

# Move selection off the edited cell, mirroring the authored workbook state
$ws.Range("M1").Select()
